# Organizing QB website 1/7
#
# Updates the "Occurrence" column (E2:E7) on Sheet1 so that question 42 is
# added to the list of occurrences, and restores the saved cell
# selection/scroll position.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: "15, 34, 38" -> "15, 34, 38. 42"
$ws.Range("E2").Value = "15, 34, 38. 42"

# Row 3: "15. 34, 38" -> "15. 34, 38, 42"
$ws.Range("E3").Value = "15. 34, 38, 42"

# Rows 4-7: "15, 34, 38" -> "15, 34, 38, 42"
$ws.Range("E4").Value = "15, 34, 38, 42"
$ws.Range("E5").Value = "15, 34, 38, 42"
$ws.Range("E6").Value = "15, 34, 38, 42"
$ws.Range("E7").Value = "15, 34, 38, 42"

# Restore the view: scrolled so column C is the leftmost visible column,
# with E7 as the active (selected) cell.
$ws.Range("E7").Select()
$excel.ActiveWindow.ScrollColumn = 3
